{"js": "// Fix a typo in the \"Nombre de proyecto\" table cell: \"TweettMonitor\" -> \"TweetMonitor\".\n// The cell holds two runs (\"Twee\" + \"ttMonitor\"); only the second run's text\n// needs the leading duplicated \"t\" removed so it reads \"tMonitor\".\nconst results = context.document.body.search(\"ttMonitor\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the matched \"ttMonitor\" run text with \"tMonitor\" (drops the extra \"t\").\n  results.items[0].insertText(\"tMonitor\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Fix a typo in the \"Nombre de proyecto\" table cell: \"TweettMonitor\" -> \"TweetMonitor\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"ttMonitor\"\n$find.Replacement.Text = \"tMonitor\"\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n\n# wdReplaceAll = 2\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2) | Out-Null\n"}
